# Daily attendance processing - 2026-01-02 06:43:59
# Swap the order of authors in the "Recorded By" column (G) for sessions
# that were recorded by both dnasr281@gmail.com and the System user:
# "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$newValue = "System, dnasr281@gmail.com"

$ws.Range("G8").Value = $newValue
$ws.Range("G9").Value = $newValue
$ws.Range("G10").Value = $newValue
$ws.Range("G12").Value = $newValue
$ws.Range("G14").Value = $newValue
$ws.Range("G15").Value = $newValue
$ws.Range("G17").Value = $newValue
$ws.Range("G18").Value = $newValue
$ws.Range("G34").Value = $newValue
$ws.Range("G35").Value = $newValue
$ws.Range("G36").Value = $newValue
$ws.Range("G38").Value = $newValue
$ws.Range("G40").Value = $newValue
$ws.Range("G41").Value = $newValue
$ws.Range("G43").Value = $newValue
$ws.Range("G44").Value = $newValue
$ws.Range("G60").Value = $newValue
$ws.Range("G61").Value = $newValue
$ws.Range("G62").Value = $newValue
$ws.Range("G64").Value = $newValue
$ws.Range("G66").Value = $newValue
$ws.Range("G67").Value = $newValue
$ws.Range("G69").Value = $newValue
$ws.Range("G70").Value = $newValue
$ws.Range("G86").Value = $newValue
$ws.Range("G87").Value = $newValue
$ws.Range("G88").Value = $newValue
$ws.Range("G90").Value = $newValue
$ws.Range("G92").Value = $newValue
$ws.Range("G93").Value = $newValue
$ws.Range("G95").Value = $newValue
$ws.Range("G96").Value = $newValue
$ws.Range("G112").Value = $newValue
$ws.Range("G113").Value = $newValue
$ws.Range("G114").Value = $newValue
$ws.Range("G116").Value = $newValue
$ws.Range("G118").Value = $newValue
$ws.Range("G119").Value = $newValue
$ws.Range("G121").Value = $newValue
$ws.Range("G122").Value = $newValue
$ws.Range("G138").Value = $newValue
$ws.Range("G139").Value = $newValue
$ws.Range("G140").Value = $newValue
$ws.Range("G142").Value = $newValue
$ws.Range("G144").Value = $newValue
$ws.Range("G145").Value = $newValue
$ws.Range("G147").Value = $newValue
$ws.Range("G148").Value = $newValue
$ws.Range("G164").Value = $newValue
$ws.Range("G167").Value = $newValue
$ws.Range("G170").Value = $newValue
$ws.Range("G174").Value = $newValue
$ws.Range("G191").Value = $newValue
$ws.Range("G194").Value = $newValue
$ws.Range("G197").Value = $newValue
$ws.Range("G201").Value = $newValue
$ws.Range("G218").Value = $newValue
$ws.Range("G221").Value = $newValue
$ws.Range("G224").Value = $newValue
$ws.Range("G228").Value = $newValue
$ws.Range("G245").Value = $newValue
$ws.Range("G248").Value = $newValue
$ws.Range("G251").Value = $newValue
$ws.Range("G255").Value = $newValue
$ws.Range("G272").Value = $newValue
$ws.Range("G275").Value = $newValue
$ws.Range("G278").Value = $newValue
$ws.Range("G282").Value = $newValue
$ws.Range("G299").Value = $newValue
$ws.Range("G302").Value = $newValue
$ws.Range("G305").Value = $newValue
$ws.Range("G309").Value = $newValue
